$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8164013333333333
$ws.Range("H2").Value = 2.449204
$ws.Range("I2").Value = 0.05618115571687973
$ws.Range("J2").Value = 0.05618115571687973
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.024021792832
$ws.Range("R2").Value = 0.216196135488
$ws.Range("S2").Value = 0.008275840567385086
$ws.Range("T2").Value = 0.008275840567385086

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8164013333333333
$ws.Range("H3").Value = 2.449204
$ws.Range("I3").Value = 0.05618115571687973
$ws.Range("J3").Value = 0.05618115571687973
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 0.1390519242973333
$ws.Range("R3").Value = 1.251467318676
$ws.Range("S3").Value = 0.04790531514949464
$ws.Range("T3").Value = 0.04790531514949464

# Row 4
$ws.Range("I4").Value = 0.8862323361798529
$ws.Range("J4").Value = 0.8862323361798529
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.378932923488
$ws.Range("R4").Value = 3.410396311392
$ws.Range("S4").Value = 0.1305476440685266
$ws.Range("T4").Value = 0.1305476440685266

# Row 5
$ws.Range("I5").Value = 0.8862323361798529
$ws.Range("J5").Value = 0.8862323361798529
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.7556846921113262
$ws.Range("T5").Value = 0.7556846921113262

# Row 6
$ws.Range("G6").Value = 0.8368233333333334
$ws.Range("I6").Value = 0.05758650810326746
$ws.Range("J6").Value = 0.05758650810326746
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.02462268976
$ws.Range("R6").Value = 0.22160420784
$ws.Range("S6").Value = 0.008482857887380244
$ws.Range("T6").Value = 0.008482857887380244

# Row 7
$ws.Range("G7").Value = 0.8368233333333334
$ws.Range("I7").Value = 0.05758650810326746
$ws.Range("J7").Value = 0.05758650810326746
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("S7").Value = 0.04910365021588721
$ws.Range("T7").Value = 0.04910365021588721
